# Plantilla_muestras.xlsx update:
#  - add new "localizacion" header columns (Congelador..Subposicion)
#  - mark the existing ID columns (A:B) and the new columns as bold+red
#  - move the active selection to M5
#  - set the page setup (paper size / orientation) for the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (P1:V1)
$ws.Range("P1").Value = "Congelador"
$ws.Range("Q1").Value = "Estante"
$ws.Range("R1").Value = "Posición del rack en el estante"
$ws.Range("S1").Value = "Rack"
$ws.Range("T1").Value = "Posición de la caja en el rack"
$ws.Range("U1").Value = "Caja"
$ws.Range("V1").Value = "Subposición"

# Highlight the ID columns and the new location columns in bold red
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Color = 255

$ws.Range("P1:V1").Font.Bold = $true
$ws.Range("P1:V1").Font.Color = 255

# Move the selection like the saved workbook (cursor left at M5)
$null = $ws.Range("M5").Select()

# Page setup: letter-ish paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
